$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 95
$lastCol = 26

# Mapping from target (after) row -> source (before) row.
# Derived from the diff: this edit is a pure re-ordering of the data rows
# (rows 2..95) combined with a +1 day bump of every "Förändrad" (column C)
# date. No other cell content actually changes.
$mapping = @{}
$mapping[2] = 2
$mapping[3] = 3
$mapping[4] = 4
$mapping[5] = 5
$mapping[6] = 6
$mapping[7] = 7
$mapping[8] = 8
$mapping[9] = 9
$mapping[10] = 10
$mapping[11] = 11
$mapping[12] = 12
$mapping[13] = 13
$mapping[14] = 14
$mapping[15] = 16
$mapping[16] = 15
$mapping[17] = 17
$mapping[18] = 18
$mapping[19] = 19
$mapping[20] = 20
$mapping[21] = 24
$mapping[22] = 21
$mapping[23] = 23
$mapping[24] = 22
$mapping[25] = 25
$mapping[26] = 26
$mapping[27] = 27
$mapping[28] = 28
$mapping[29] = 29
$mapping[30] = 30
$mapping[31] = 31
$mapping[32] = 32
$mapping[33] = 33
$mapping[34] = 34
$mapping[35] = 35
$mapping[36] = 36
$mapping[37] = 37
$mapping[38] = 38
$mapping[39] = 39
$mapping[40] = 40
$mapping[41] = 41
$mapping[42] = 42
$mapping[43] = 86
$mapping[44] = 89
$mapping[45] = 72
$mapping[46] = 74
$mapping[47] = 64
$mapping[48] = 49
$mapping[49] = 81
$mapping[50] = 55
$mapping[51] = 56
$mapping[52] = 62
$mapping[53] = 59
$mapping[54] = 60
$mapping[55] = 61
$mapping[56] = 52
$mapping[57] = 48
$mapping[58] = 50
$mapping[59] = 75
$mapping[60] = 57
$mapping[61] = 67
$mapping[62] = 53
$mapping[63] = 43
$mapping[64] = 63
$mapping[65] = 87
$mapping[66] = 85
$mapping[67] = 71
$mapping[68] = 70
$mapping[69] = 73
$mapping[70] = 77
$mapping[71] = 90
$mapping[72] = 79
$mapping[73] = 54
$mapping[74] = 84
$mapping[75] = 65
$mapping[76] = 78
$mapping[77] = 66
$mapping[78] = 45
$mapping[79] = 88
$mapping[80] = 46
$mapping[81] = 83
$mapping[82] = 80
$mapping[83] = 91
$mapping[84] = 69
$mapping[85] = 44
$mapping[86] = 58
$mapping[87] = 82
$mapping[88] = 68
$mapping[89] = 76
$mapping[90] = 92
$mapping[91] = 95
$mapping[92] = 93
$mapping[93] = 94
$mapping[94] = 47
$mapping[95] = 51

# 1) Snapshot every cell in rows 2..95 (columns A..Z) before making any
#    changes, so that overwriting a row does not destroy data that is
#    still needed as the source for another target row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $h = @{}
        if ($cell.HasFormula) {
            $h["IsFormula"] = $true
            $h["Data"] = $cell.Formula
        } else {
            $h["IsFormula"] = $false
            $h["Data"] = $cell.Value2
        }
        $rowData += $h
    }
    $snapshot[$r] = $rowData
}

# 2) Write the snapshot data back out in the new (permuted) row order.
for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $mapping[$targetRow]
    $rowData = $snapshot[$sourceRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $src = $rowData[$c - 1]
        $cell = $ws.Cells.Item($targetRow, $c)
        if ($src["Data"] -eq $null) {
            $cell.ClearContents()
        } elseif ($src["IsFormula"]) {
            $cell.Formula = $src["Data"]
        } else {
            $cell.Value2 = $src["Data"]
        }
    }
}

# 3) Bump the "Förändrad" (column C) date by one day for every data row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
